$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel alignment constants
$xlLeft = -4131
$xlCenter = -4108

# --- New "navigation" test case descriptions added to column H (rows 26-51) ---
# Column H mirrors the wrap/left/center style already used by column D (style index 6).
# NOTE: cells are populated in this specific order (not strict row order) so that the
# workbook shared-string table is rebuilt in the same sequence as the source edit.
$cell = $ws.Cells.Item(26, 8)
$cell.Value = "the navigation bar appears in the page in it's position"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(28, 8)
$cell.Value = "the search bar appears in the page in it's position"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(27, 8)
$cell.Value = "the logo  appears in the page in it's position"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(29, 8)
$cell.Value = "the menu items appear in`n their order"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(30, 8)
$cell.Value = "the `"home `" button direct the user to the home page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(31, 8)
$cell.Value = "the `" Software `" button direct the user to the Software page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(32, 8)
$cell.Value = "the `" Networking `" button direct the user to the Networking page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(33, 8)
$cell.Value = "the `" Embeded Syaytems `" button direct the user to the Embeded Syaytems page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(34, 8)
$cell.Value = "the `" Biotechnology `" button direct the user to the Biotechnology page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(35, 8)
$cell.Value = "the`" My Profile `" button direct the user to the his/her Profile page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(36, 8)
$cell.Value = "the `" ADD atricle `" button direct the user to the ADD atricle page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(37, 8)
$cell.Value = "the `" Notification `" button direct the user to the Notification page"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(38, 8)
$cell.Value = "search result appear with `nauther first name"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(39, 8)
$cell.Value = "search result appear with `nauther last name"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(40, 8)
$cell.Value = "search result appear with `nauther first and last name"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(41, 8)
$cell.Value = "search result appear with `narticle tittle"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(42, 8)
$cell.Value = "search results apper under`nthe search bar"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(43, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(44, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(45, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(46, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(47, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(48, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(49, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(50, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

$cell = $ws.Cells.Item(51, 8)
$cell.Value = "the response time is as `nrequired"
$cell.WrapText = $true
$cell.HorizontalAlignment = $xlLeft
$cell.VerticalAlignment = $xlCenter

# --- Row heights grow to fit the newly wrapped column-H text ---
$ws.Rows.Item(29).RowHeight = 28.8
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 43.2
$ws.Rows.Item(33).RowHeight = 57.6
$ws.Rows.Item(34).RowHeight = 43.2
$ws.Rows.Item(35).RowHeight = 43.2
$ws.Rows.Item(36).RowHeight = 43.2
$ws.Rows.Item(37).RowHeight = 43.2
$ws.Rows.Item(41).RowHeight = 35.4

# --- Update the active selection to reflect where the edits were made ---
$ws.Range("H51").Select()
